$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C2").Value2 = "Final dilution and order of magnitude specified"
$v1 = $ws2.Range("C2").Value2
Write-Host "C2 Value2 after set: $v1"

$ws2.Range("I1").Value2 = 20230209
$v2 = $ws2.Range("I1").Value2
Write-Host "I1 Value2 after set: $v2"
